# MCH337 collections update
# Adds two new "Series" records (MCH337-1, MCH337-2) to the finding-aid
# sheet, matching the formatting already used for the sibling MCH
# workbooks, then re-establishes the frozen header row / active
# selection the way Excel leaves them after such an edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- formatting -----------------------------------------------------
# Column B (alternativeIdentifiers) is intentionally left untouched/blank
# for these rows, so format the identifier column and the rest of the
# record separately.
$idCol = $ws.Range("A2:A3")
$idCol.Font.Name = "Calibri"
$idCol.Font.ThemeColor = 1

$rest = $ws.Range("C2:H3")
$rest.Font.Name = "Calibri"
$rest.Font.ThemeColor = 1

# --- row 2: MCH337-1 --------------------------------------------------
$ws.Range("A2").Value = "MCH337-1"
$ws.Range("C2").Value = "OSCAR MPETHA TRIAL"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 33H | GRAP COUNT NUMER: NONE"

# --- row 3: MCH337-2 --------------------------------------------------
$ws.Range("A3").Value = "MCH337-2"
$ws.Range("C3").Value = "OSCAR MPETHA TRIAL"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 33H | GRAP COUNT NUMER: NONE"

# --- view state: keep header frozen, select the rows just entered ----
$ws.Range("A2:I3").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
